$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 5490
$ws1.Range("F8").Value = 901
$ws1.Range("F10").Value = 2446
$ws1.Range("F11").Value = 81
$ws1.Range("F12").Value = 72
$ws1.Range("F13").Value = 66
$ws1.Range("F14").Value = 2294
$ws1.Range("F15").Value = 186

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 5490
$ws4.Range("F10").Value = 901
$ws4.Range("F12").Value = 2446
$ws4.Range("F13").Value = 81
$ws4.Range("F14").Value = 72
$ws4.Range("F16").Value = 66
$ws4.Range("F17").Value = 2294
$ws4.Range("F18").Value = 186
